$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D to make room for the new quarter (shifts D:K -> E:L)
$ws.Columns("D").Insert()

# Copy number formatting (date / number styles) from the old column (now E) into the new column D
$ws.Range("E7:E102").Copy() | Out-Null
$ws.Range("D7:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the new column D with the new quarter values
$ws.Range("D7").Value = 43407
$ws.Range("D8").Value = 1857300
$ws.Range("D9").Value = 1333700
$ws.Range("D10").Value = 523600
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 1804400
$ws.Range("D18").Value = 52900
$ws.Range("D20").Value = -100
$ws.Range("D21").Value = 111500
$ws.Range("D22").Value = 2600
$ws.Range("D23").Value = 50200
$ws.Range("D24").Value = 12400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 37800
$ws.Range("D27").Value = 37800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 100
$ws.Range("D33").Value = 37800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 37800
$ws.Range("D38").Value = 43407
$ws.Range("D41").Value = 92100
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 68000
$ws.Range("D44").Value = 2196800
$ws.Range("D45").Value = 138500
$ws.Range("D46").Value = 2495300
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 1578300
$ws.Range("D49").Value = 382200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 127900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 4583800
$ws.Range("D57").Value = 1028200
$ws.Range("D58").Value = 5300
$ws.Range("D59").Value = 520600
$ws.Range("D60").Value = 1554100
$ws.Range("D61").Value = 438400
$ws.Range("D62").Value = 744200
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 2736700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 2374300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1847000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43407
$ws.Range("D81").Value = 37800
$ws.Range("D83").Value = 58700
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -155400
$ws.Range("D91").Value = -38800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -38800
$ws.Range("D96").Value = -22100
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 162000
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -32200

# Restated historical values for "Capital Expenditures" row 91 (columns G:J no longer a pure shift)
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = -2300
